# Adds a new "WS 2020/21" semester block (5 country rows + 1 total row)
# to the bottom of the data table on Tabelle1, pushing the footnote /
# source block down by 6 rows, and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 6 fresh rows right before the old trailing spacer row (101) ---
$ws.Rows("101:106").Insert()

# --- 2. Stash the (now shifted) hyperlink cell's original formatting so we
#        can restore it after re-creating the hyperlink (Hyperlinks.Add
#        resets the cell style to the generic built-in "Link" look). ---
$ws.Range("B115").Copy()
$stash = $ws.Range("Z1")
$stash.PasteSpecial(-4122)

# --- 3. Copy the format of the last complete semester block (rows 95:100,
#        "WS 2019/20") down onto the 6 new rows so borders / fonts / number
#        formats match the rest of the table. ---
$ws.Range("B95:L100").Copy()
$ws.Range("B101:L106").PasteSpecial(-4122)
$ws.Rows("101:105").RowHeight = 8.25
$ws.Rows("106").RowHeight = 16.5

# --- 4. Fill in the new "WS 2020/21" values ---
$ws.Range("B101").Value = "China"
$ws.Range("C101").Value = "WS 2020/21"
$ws.Range("D101").Value = 4286
$ws.Range("E101").Value = 2523
$ws.Range("F101").Value = 1763
$ws.Range("G101").Value = 4157
$ws.Range("H101").Value = 2451
$ws.Range("I101").Value = 1706
$ws.Range("J101").Value = 129
$ws.Range("K101").Value = 72
$ws.Range("L101").Value = 57

$ws.Range("B102").Value = "Türkei"
$ws.Range("C102").Value = "WS 2020/21"
$ws.Range("D102").Value = 1919
$ws.Range("E102").Value = 979
$ws.Range("F102").Value = 940
$ws.Range("G102").Value = 518
$ws.Range("H102").Value = 310
$ws.Range("I102").Value = 208
$ws.Range("J102").Value = 1401
$ws.Range("K102").Value = 669
$ws.Range("L102").Value = 732

$ws.Range("B103").Value = "Syrien Arab. Republik"
$ws.Range("C103").Value = "WS 2020/21"
$ws.Range("D103").Value = 1751
$ws.Range("E103").Value = 1336
$ws.Range("F103").Value = 415
$ws.Range("G103").Value = 1528
$ws.Range("H103").Value = 1199
$ws.Range("I103").Value = 329
$ws.Range("J103").Value = 223
$ws.Range("K103").Value = 137
$ws.Range("L103").Value = 86

$ws.Range("B104").Value = "Indien"
$ws.Range("C104").Value = "WS 2020/21"
$ws.Range("D104").Value = 1502
$ws.Range("E104").Value = 1080
$ws.Range("F104").Value = 422
$ws.Range("G104").Value = 1493
$ws.Range("H104").Value = 1074
$ws.Range("I104").Value = 419
$ws.Range("J104").Value = 9
$ws.Range("K104").Value = 6
$ws.Range("L104").Value = 3

$ws.Range("B105").Value = "Iran Islamische Republik"
$ws.Range("C105").Value = "WS 2020/21"
$ws.Range("D105").Value = 1124
$ws.Range("E105").Value = 586
$ws.Range("F105").Value = 538
$ws.Range("G105").Value = 1062
$ws.Range("H105").Value = 554
$ws.Range("I105").Value = 508
$ws.Range("J105").Value = 62
$ws.Range("K105").Value = 32
$ws.Range("L105").Value = 30

$ws.Range("B106").Value = "Staaten insgesamt"
$ws.Range("C106").Value = "WS 2020/21"
$ws.Range("D106").Value = 25111
$ws.Range("E106").Value = 14481
$ws.Range("F106").Value = 10630
$ws.Range("G106").Value = 20647
$ws.Range("H106").Value = 12329
$ws.Range("I106").Value = 8318
$ws.Range("J106").Value = 4464
$ws.Range("K106").Value = 2152
$ws.Range("L106").Value = 2312

# --- 5. Re-point the "source" hyperlink at its new location (B115) and
#        restore its original look. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B115"), "https://www.integrationsmonitoring.niedersachsen.de/")
$stash.Copy()
$ws.Range("B115").PasteSpecial(-4122)
$stash.Clear()

# --- 6. Update the active selection to match the author's final cursor
#        position. ---
$ws.Range("B106").Select()
